$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header row: "_old" -> "_FV2304", "_new" -> "_FV2310"
# (column "diff" in K stays unchanged)
$headers = @(
    "Segmentname_FV2304","Segmentgruppe_FV2304","Segment_FV2304","Datenelement_FV2304","Segment ID_FV2304",
    "Code_FV2304","Qualifier_FV2304","Beschreibung_FV2304","Bedingungsausdruck_FV2304","Bedingung_FV2304",
    "diff",
    "Segmentname_FV2310","Segmentgruppe_FV2310","Segment_FV2310","Datenelement_FV2310","Segment ID_FV2310",
    "Code_FV2310","Qualifier_FV2310","Beschreibung_FV2310","Bedingungsausdruck_FV2310","Bedingung_FV2310"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Turn the used range into an Excel Table ("Table1") so the headers drive
# the table's column names / autofilter.
$dataRange = $ws.Range("A1:U58")
$tbl = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# Freeze the header row (row 1) so it stays visible while scrolling.
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Output "done"
